$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of rows whose match data (everything except id/Div/Div Original Name/Date)
# needs to be swapped between the two rows.
$pairs = @(
    @(12, 13),
    @(19, 20),
    @(50, 51),
    @(100, 101),
    @(173, 174)
)

# Columns to swap: B (2) and F..AC (6..29). Columns A, C, D, E stay untouched.
$cols = @(2) + @(6..29)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($c in $cols) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
